$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 9.006948722776087
$ws.Range("G2").Value = 8.806415168460042
$ws.Range("H2").Value = 9.208468791592765
$ws.Range("I2").Value = 0.009749600947258831
$ws.Range("J2").Value = 0.008728229721266247
$ws.Range("K2").Value = 0.01091254621706741
$ws.Range("L2").Value = 0.004702502597597281
$ws.Range("M2").Value = 0.004457556515041075
$ws.Range("N2").Value = 0.004966887922525317

# Row 3
$ws.Range("F3").Value = 0.2773227466811037
$ws.Range("G3").Value = 0.276475661748702
$ws.Range("H3").Value = 0.2782010286090696
$ws.Range("I3").Value = 0.2494150788188885
$ws.Range("J3").Value = 0.2486165673564611
$ws.Range("K3").Value = 0.2502412475741974
$ws.Range("L3").Value = 0.2747274822470684
$ws.Range("M3").Value = 0.273884035153016
$ws.Range("N3").Value = 0.275602624241971

# Row 4
$ws.Range("F4").Value = 9.284271469457192
$ws.Range("G4").Value = 9.082890830208743
$ws.Range("H4").Value = 9.486669820201833
$ws.Range("I4").Value = 0.2591646797661473
$ws.Range("J4").Value = 0.2573447970777274
$ws.Range("K4").Value = 0.2611537937912648
$ws.Range("L4").Value = 0.2794299848446657
$ws.Range("M4").Value = 0.2783415916680571
$ws.Range("N4").Value = 0.2805695121644963
